# Applies the "fixing links to move off statistics.byuimath" edit.
# The change set (per the source diff) is a handful of small, surgical
# OOXML tweaks scattered across word/document.xml, word/numbering.xml and
# word/styles.xml:
#   1. w:bCs added next to the existing w:b on the "Please note..." run.
#   2. w:tblLook gains a w:val="0020" attribute.
#   3. The hard-coded first-row formatting (w:cnfStyle on the row, plus
#      w:tcBorders/w:vAlign on each of its three cells) is removed from
#      the table markup in the body...
#   4. ...and reappears as a w:tblStylePr (type="firstRow") conditional
#      format inside the "Table" table style definition.
#   5. Three m:r runs in the oMath expressions gain an m:rPr/m:sty
#      (val="p") so they render with "plain" (non-italic) style.
#   6. The bullet-list numbering definition (abstractNum 990) drops its
#      w:tabs/w:tab children and widens each level's hanging indent by
#      240 twips (480->720, 1200->1440, ... 6240->6480).
#
# Because several of these land in parts that aren't reachable through
# the high-level Word object model (conditional table-style formatting,
# numbering level tab stops), we round-trip the whole package through
# Document.Content.WordOpenXML, a flat-OPC rendering of every part
# (document.xml, styles.xml, numbering.xml, ...), edit it as text, and
# write it back.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1) Bold run also gets bCs (bold-complex-script) turned on.
# ---------------------------------------------------------------------
$old = '<w:rPr><w:b/></w:rPr><w:t>Please note that the steps show rounded numbers,'
$new = '<w:rPr><w:b/><w:bCs/></w:rPr><w:t>Please note that the steps show rounded numbers,'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: bCs run" }
$xml = $xml.Replace($old, $new)

# ---------------------------------------------------------------------
# 2) tblLook picks up an explicit w:val bitmask.
# ---------------------------------------------------------------------
$old = '<w:tblLook w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>'
$new = '<w:tblLook w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0" w:val="0020"/>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: tblLook" }
$xml = $xml.Replace($old, $new)

# ---------------------------------------------------------------------
# 3) Drop the direct first-row formatting - it moves into the style.
# ---------------------------------------------------------------------
$old = '<w:trPr><w:cnfStyle w:firstRow="1" /></w:trPr>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: cnfStyle" }
$xml = $xml.Replace($old, '')

$old = '<w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: header tcPr" }
$xml = $xml.Replace($old, '')

# ---------------------------------------------------------------------
# 4) ...and add it back as conditional "firstRow" table-style formatting.
# ---------------------------------------------------------------------
$old = '<w:style w:type="table" w:default="1" w:styleId="Table"><w:name w:val="Table"/><w:basedOn w:val="TableNormal"/><w:semiHidden/><w:unhideWhenUsed/><w:qFormat/><w:tblPr>' + "`n      " + '<w:tblInd w:w="0" w:type="dxa" />' + "`n      " + '<w:tblCellMar>' + "`n        " + '<w:top w:w="0" w:type="dxa" />' + "`n        " + '<w:left w:w="108" w:type="dxa" />' + "`n        " + '<w:bottom w:w="0" w:type="dxa" />' + "`n        " + '<w:right w:w="108" w:type="dxa" />' + "`n      " + '</w:tblCellMar>' + "`n    " + '</w:tblPr></w:style>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: Table style" }
$new = $old.Replace('</w:tblPr></w:style>', '</w:tblPr><w:tblStylePr w:type="firstRow"><w:tblPr><w:jc w:val="left"/><w:tblInd w:w="0" w:type="dxa"/></w:tblPr><w:trPr><w:jc w:val="left"/></w:trPr><w:tcPr><w:vAlign w:val="bottom"/><w:tcBorders><w:bottom w:val="single"/></w:tcBorders></w:tcPr></w:tblStylePr></w:style>')
$xml = $xml.Replace($old, $new)

# ---------------------------------------------------------------------
# 5) oMath runs: give the bare punctuation runs the "plain" math style.
# ---------------------------------------------------------------------
$old = '<m:r><m:t>' + [char]0x2212 + '</m:t></m:r>'
$new = '<m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>' + [char]0x2212 + '</m:t></m:r>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: minus run" }
$xml = $xml.Replace($old, $new)

$old = '<m:r><m:t>(</m:t></m:r>'
$new = '<m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>(</m:t></m:r>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: open paren run" }
$xml = $xml.Replace($old, $new)

$old = '<m:sSup><m:e><m:r><m:t>)</m:t></m:r></m:e>'
$new = '<m:sSup><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>)</m:t></m:r></m:e>'
if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: close paren run" }
$xml = $xml.Replace($old, $new)

# ---------------------------------------------------------------------
# 6) Bullet list levels: drop the w:tabs, widen the hanging indent.
# ---------------------------------------------------------------------
$pairs = @(
    @("0", "480", "720"),
    @("720", "1200", "1440"),
    @("1440", "1920", "2160"),
    @("2160", "2640", "2880"),
    @("2880", "3360", "3600"),
    @("3600", "4080", "4320"),
    @("4320", "4800", "5040"),
    @("5040", "5520", "5760"),
    @("5760", "6240", "6480")
)
foreach ($p in $pairs) {
    $tabPos = $p[0]
    $oldLeft = $p[1]
    $newLeft = $p[2]
    $old = '<w:pPr><w:tabs><w:tab w:val="num" w:pos="' + $tabPos + '" /></w:tabs><w:ind w:left="' + $oldLeft + '" w:hanging="480" /></w:pPr>'
    $new = '<w:pPr><w:ind w:left="' + $newLeft + '" w:hanging="480"/></w:pPr>'
    if ($xml.IndexOf($old) -lt 0) { throw "pattern not found: numbering level tabPos=$tabPos" }
    $xml = $xml.Replace($old, $new)
}

$d.Content.WordOpenXML = $xml
Write-Output "applied edits"
